$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 holds literal text that happens to start with "=" (not a real formula).
# Update the text (addAll -> addAll1) while keeping it as literal text by
# prefixing with an apostrophe; Excel marks such cells with quotePrefix="1"
# in the resulting style.
$ws.Range("B8").Value = "'= addAll1(null, null); ""Hello"";"

# Move the active selection to E7.
$ws.Range("E7").Select()
